# Updates cryptos list data in columns B-E for rows 2-51 per the diff.
# Column D cells that are pure numeric-looking strings (e.g. "603.55") are
# written with a leading apostrophe so Excel keeps them as text (matching
# the original inlineStr/text cell type) instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.881.50'
$ws.Range('E2').Value = '  -0.50%  '

$ws.Range('D3').Value = '3.532.26'
$ws.Range('E3').Value = '  -0.48%  '

$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.23%  '

$ws.Range('D5').Value = '''603.55'
$ws.Range('E5').Value = '  -1.91%  '

$ws.Range('D6').Value = '''195.52'
$ws.Range('E6').Value = '  +3.40%  '

$ws.Range('D7').Value = '''0.633'
$ws.Range('E7').Value = '  -0.72%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('D9').Value = '''0.205'
$ws.Range('E9').Value = '  -5.10%  '

$ws.Range('D10').Value = '''0.652'
$ws.Range('E10').Value = '  -2.01%  '

$ws.Range('E11').Value = '  -0.47%  '

$ws.Range('D12').Value = '''0.0000301'
$ws.Range('E12').Value = '  -2.58%  '

$ws.Range('D13').Value = '''9.52'
$ws.Range('E13').Value = '  -2.42%  '

$ws.Range('D14').Value = '4.094.04'
$ws.Range('E14').Value = '  -0.58%  '

$ws.Range('D15').Value = '''611.99'
$ws.Range('E15').Value = '  -0.39%  '

$ws.Range('D16').Value = '''12.85'
$ws.Range('E16').Value = '  -0.22%  '

$ws.Range('E17').Value = '  -0.16%  '

$ws.Range('D18').Value = '70.075.11'
$ws.Range('E18').Value = '  -0.39%  '

$ws.Range('D19').Value = '3.542.82'
$ws.Range('E19').Value = '  -1.33%  '

$ws.Range('E20').Value = '  +0.25%  '

$ws.Range('D21').Value = '''0.995'
$ws.Range('E21').Value = '  -0.80%  '

$ws.Range('D22').Value = '''18.12'
$ws.Range('E22').Value = '  +2.62%  '

$ws.Range('D23').Value = '''5.31'
$ws.Range('E23').Value = '  +3.36%  '

$ws.Range('D24').Value = '''102.74'
$ws.Range('E24').Value = '  -2.93%  '

$ws.Range('E25').Value = '  -2.17%  '

$ws.Range('D26').Value = '''3.14'
$ws.Range('E26').Value = '  +2.93%  '

$ws.Range('D27').Value = '''10.91'
$ws.Range('E27').Value = '  -0.76%  '

$ws.Range('D28').Value = '''9.58'
$ws.Range('E28').Value = '  -5.81%  '

$ws.Range('E29').Value = '  -2.92%  '

$ws.Range('D30').Value = '''7.12'
$ws.Range('E30').Value = '  -0.10%  '

$ws.Range('D31').Value = '''4.31'
$ws.Range('E31').Value = '  +14.87%  '

$ws.Range('E32').Value = '  +0.72%  '

$ws.Range('E33').Value = '  -2.11%  '

$ws.Range('D34').Value = '''63.26'
$ws.Range('E34').Value = '  -1.66%  '

$ws.Range('D35').Value = '0.0₃0859'
$ws.Range('E35').Value = '  +9.81%  '

$ws.Range('D36').Value = '3.725.24'
$ws.Range('E36').Value = '  +5.05%  '

$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '''3.66'
$ws.Range('E38').Value = '  +2.40%  '

$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '''3.03'
$ws.Range('E39').Value = '  -4.19%  '

$ws.Range('E40').Value = '  -2.20%  '

$ws.Range('D41').Value = '''36.56'
$ws.Range('E41').Value = '  -1.95%  '

$ws.Range('D42').Value = '''491.86'
$ws.Range('E42').Value = '  -8.84%  '

$ws.Range('E43').Value = '  -4.85%  '

$ws.Range('E44').Value = '  -2.36%  '

$ws.Range('E45').Value = '  -1.92%  '

$ws.Range('D46').Value = '''2.84'
$ws.Range('E46').Value = '  -4.56%  '

$ws.Range('D47').Value = '''3.30'
$ws.Range('E47').Value = '  -1.95%  '

$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.24%  '

$ws.Range('E49').Value = '  -4.21%  '

$ws.Range('E50').Value = '  +4.34%  '

$ws.Range('D51').Value = '''130.87'
$ws.Range('E51').Value = '  -1.36%  '
